$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing the duplicate/extraneous message
# "¿Cuál es tu gracias?" (classified as "nombre"), which was row 4.
# All subsequent rows shift up by one.
$ws.Rows.Item(4).Delete() | Out-Null

# Update the active selection to match the saved view state.
$ws.Range("B7").Select() | Out-Null
